$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the "取得日時" (fetched at) timestamp on the existing rows to the
# new run time.
$ws.Range("A2").Value = "2025-12-27 12:35:43"
$ws.Range("A3").Value = "2025-12-27 12:35:43"

# Append the two newly scraped listings as rows 4 and 5.
$ws.Range("A4").Value = "2025-12-27 12:35:43"
$ws.Range("B4").Value = "【募集】Web予約フロー自動化ツールの設計・開発をお任せします"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5462249"
$ws.Range("G4").Value = 213
$ws.Range("H4").Value = "◆ツール,開発"

$ws.Range("A5").Value = "2025-12-27 12:35:43"
$ws.Range("B5").Value = "【急募】ReactでLine風会話履歴表示コンポーネント作成依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5462198"
$ws.Range("G5").Value = 128
$ws.Range("H5").Value = "🔥React"

# Wire up the hyperlinks for the new URL cells, matching the style used by
# the existing hyperlink cells.
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5462249")
$ws.Range("F4").Style = $ws.Range("F2").Style

$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5462198")
$ws.Range("F5").Style = $ws.Range("F2").Style
